$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Range("G10").Style = "Good"
$ws.Range("G10").Value = "√"
Write-Host "done"
